# Consolidate the "men_cloths" row into the Nike/Footwear summary row:
#  - delete row 3 (Nike, men_cloths, cloths, 200, 6, 1200)
#  - roll its totals into row 2 (Nike, Footwear, ...)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()

$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2000
